# Auto-generated: applies 106 cell-value updates to sheet1 (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.420.18"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.989.02"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'385.03"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'102.56"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "'0.542"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D9").Value = "'0.595"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "'36.83"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.0845"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "3.459.69"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "'18.24"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'7.52"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "2.989.48"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "'1.00"
$ws.Range("E17").Value = "  +6.91%  "
$ws.Range("D18").Value = "51.374.85"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'3.28"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("D20").Value = "'7.42"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "'12.91"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'68.92"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "'262.21"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'2.92"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("D26").Value = "'8.22"
$ws.Range("E26").Value = "  +13.54%  "
$ws.Range("D27").Value = "'7.47"
$ws.Range("E27").Value = "  +7.88%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'4.14"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'0.168"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.113"
$ws.Range("E30").Value = "  +10.58%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'25.91"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'9.86"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "'34.58"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "'50.96"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").Value = "'2.07"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").Value = "'0.0451"
$ws.Range("E37").Value = "  +6.04%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'17.18"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'2.61"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").Value = "'1.82"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").Value = "'122.41"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "'21.68"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").Value = "'2.08"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "'0.274"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.36"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'3.28"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.033.00"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'0.0334"
$ws.Range("E51").Value = "  +0.80%  "
